$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for each coin row with the
# latest scraped values. Some new Price values are purely numeric-looking
# strings (e.g. "0.999", "6.70"); these cells are explicitly formatted as Text
# first so Excel keeps them as literal strings (matching the source data)
# instead of silently converting them to numbers and losing formatting such
# as trailing zeros.

$ws.Range("D2").Value = "57.185.42"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "3.037.61"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.19"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.58"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.19"
$ws.Range("E9").Value = "  -5.42%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.376"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").Value = "3.565.29"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.96"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "57.045.85"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "3.036.11"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("E19").Value = "  +4.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.11"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.38"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.506"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.25"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").Value = "3.167.48"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "0.0₃0885"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.70"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.73"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.44"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0676"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.660"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").Value = "2.199.34"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.09"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.944"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0865"
$ws.Range("E51").Value = "  -3.23%  "

Write-Host "Updated cryptos list"
